$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 27.295786
$ws.Range("H2").Value = 81.88735800000001
$ws.Range("I2").Value = 0.7928847908394133
$ws.Range("J2").Value = 0.7928847908394133
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("Q2").Value = 2.249509414427334
$ws.Range("R2").Value = 20.245584729846
$ws.Range("S2").Value = 0.03747970856650421
$ws.Range("T2").Value = 0.03747970856650421

# Row 3
$ws.Range("G3").Value = 27.295786
$ws.Range("H3").Value = 81.88735800000001
$ws.Range("I3").Value = 0.7928847908394133
$ws.Range("J3").Value = 0.7928847908394133
$ws.Range("Q3").Value = 45.33895564486401
$ws.Range("R3").Value = 408.050600803776
$ws.Range("S3").Value = 0.7554050822729091
$ws.Range("T3").Value = 0.755405082272909

# Row 4
$ws.Range("I4").Value = 0.03666380186764524
$ws.Range("J4").Value = 0.03666380186764524
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("R4").Value = 0.9361764985350001
$ws.Range("S4").Value = 0.001733099972172008
$ws.Range("T4").Value = 0.001733099972172009

# Row 5
$ws.Range("I5").Value = 0.03666380186764524
$ws.Range("J5").Value = 0.03666380186764524
$ws.Range("S5").Value = 0.03493070189547323
$ws.Range("T5").Value = 0.03493070189547323

# Row 6
$ws.Range("G6").Value = 3.215569666666667
$ws.Range("H6").Value = 9.646709
$ws.Range("I6").Value = 0.0934054905978733
$ws.Range("J6").Value = 0.0934054905978733
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("Q6").Value = 0.2650025992258889
$ws.Range("R6").Value = 2.385023393033
$ws.Range("S6").Value = 0.004415282783282288
$ws.Range("T6").Value = 0.004415282783282288

# Row 7
$ws.Range("G7").Value = 3.215569666666667
$ws.Range("H7").Value = 9.646709
$ws.Range("I7").Value = 0.0934054905978733
$ws.Range("J7").Value = 0.0934054905978733
$ws.Range("Q7").Value = 5.341138390005334
$ws.Range("R7").Value = 48.070245510048
$ws.Range("S7").Value = 0.08899020781459102
$ws.Range("T7").Value = 0.08899020781459101

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2487956666666666
$ws.Range("H8").Value = 0.7463869999999999
$ws.Range("I8").Value = 0.007226987349869769
$ws.Range("J8").Value = 0.007226987349869769
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("Q8").Value = 0.02050383141322222
$ws.Range("R8").Value = 0.184534482719
$ws.Range("S8").Value = 0.0003416200976691343
$ws.Range("T8").Value = 0.0003416200976691343

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2487956666666666
$ws.Range("H9").Value = 0.7463869999999999
$ws.Range("I9").Value = 0.007226987349869769
$ws.Range("J9").Value = 0.007226987349869769
$ws.Range("Q9").Value = 0.4132555734293333
$ws.Range("R9").Value = 3.719300160864
$ws.Range("S9").Value = 0.006885367252200635
$ws.Range("T9").Value = 0.006885367252200634

# Row 10
$ws.Range("G10").Value = 2.403580666666667
$ws.Range("H10").Value = 7.210742
$ws.Range("I10").Value = 0.06981892934519847
$ws.Range("J10").Value = 0.06981892934519847
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("Q10").Value = 0.1980846910948889
$ws.Range("R10").Value = 1.782762219854
$ws.Range("S10").Value = 0.003300344708987334
$ws.Range("T10").Value = 0.003300344708987334

# Row 11
$ws.Range("G11").Value = 2.403580666666667
$ws.Range("H11").Value = 7.210742
$ws.Range("I11").Value = 0.06981892934519847
$ws.Range("J11").Value = 0.06981892934519847
$ws.Range("Q11").Value = 3.992405173269333
$ws.Range("R11").Value = 35.931646559424
$ws.Range("S11").Value = 0.06651858463621114
$ws.Range("T11").Value = 0.06651858463621113
